# Refresh scraped career-stat rows (new game played, new box-score totals)
# across the four stat tables: per_game, per_minute, per_poss, advanced.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("per_game")
$ws.Range("G11").Value = 67
$ws.Range("I11").Value = 17.9
$ws.Range("K11").Value = 5.8
$ws.Range("L11").Value = 0.382
$ws.Range("M11").Value = 1.7
$ws.Range("N11").Value = 4.6
$ws.Range("R11").Value = 0.458
$ws.Range("S11").Value = 0.524
$ws.Range("V11").Value = 0.759
$ws.Range("Z11").Value = 2.6
$ws.Range("AD11").Value = 1.1
$ws.Range("AE11").Value = 7.4
$ws.Range("G13").Value = 14
$ws.Range("I13").Value = 27.9
$ws.Range("K13").Value = 11.4
$ws.Range("L13").Value = 0.403
$ws.Range("N13").Value = 8.4
$ws.Range("O13").Value = 0.385
$ws.Range("R13").Value = 0.452
$ws.Range("S13").Value = 0.544
$ws.Range("U13").Value = 4.1
$ws.Range("V13").Value = 0.772
$ws.Range("X13").Value = 2.2
$ws.Range("Y13").Value = 2.7
$ws.Range("Z13").Value = 4.4
$ws.Range("AA13").Value = 0.6
$ws.Range("AC13").Value = 1.4
$ws.Range("G20").Value = 14
$ws.Range("I20").Value = 27.9
$ws.Range("K20").Value = 11.4
$ws.Range("L20").Value = 0.403
$ws.Range("N20").Value = 8.4
$ws.Range("O20").Value = 0.385
$ws.Range("R20").Value = 0.452
$ws.Range("S20").Value = 0.544
$ws.Range("U20").Value = 4.1
$ws.Range("V20").Value = 0.772
$ws.Range("X20").Value = 2.2
$ws.Range("Y20").Value = 2.7
$ws.Range("Z20").Value = 4.4
$ws.Range("AA20").Value = 0.6
$ws.Range("AC20").Value = 1.4
$ws.Range("G22").Value = 307
$ws.Range("U22").Value = 2.3
$ws.Range("V22").Value = 0.8149999999999999
$ws.Range("AE22").Value = 11.7
$ws.Range("U24").Value = -0.9999999999999998
$ws.Range("V24").Value = 0.06000000000000005
$ws.Range("AE24").Value = -7.699999999999999

$ws = $wb.Worksheets.Item("per_minute")
$ws.Range("G11").Value = 67
$ws.Range("I11").Value = 1202
$ws.Range("J11").Value = 4.5
$ws.Range("K11").Value = 11.7
$ws.Range("L11").Value = 0.382
$ws.Range("N11").Value = 9.199999999999999
$ws.Range("Q11").Value = 2.5
$ws.Range("R11").Value = 0.458
$ws.Range("T11").Value = 3.5
$ws.Range("U11").Value = 0.759
$ws.Range("W11").Value = 2.7
$ws.Range("Y11").Value = 5.3
$ws.Range("AA11").Value = 0.4
$ws.Range("AB11").Value = 1.4
$ws.Range("AD11").Value = 14.9
$ws.Range("G13").Value = 14
$ws.Range("I13").Value = 391
$ws.Range("K13").Value = 14.6
$ws.Range("L13").Value = 0.403
$ws.Range("N13").Value = 10.8
$ws.Range("O13").Value = 0.385
$ws.Range("P13").Value = 1.7
$ws.Range("Q13").Value = 3.9
$ws.Range("R13").Value = 0.452
$ws.Range("S13").Value = 4.1
$ws.Range("T13").Value = 5.2
$ws.Range("U13").Value = 0.772
$ws.Range("Y13").Value = 5.6
$ws.Range("Z13").Value = 0.8
$ws.Range("AD13").Value = 20
$ws.Range("G20").Value = 14
$ws.Range("I20").Value = 391
$ws.Range("K20").Value = 14.6
$ws.Range("L20").Value = 0.403
$ws.Range("N20").Value = 10.8
$ws.Range("O20").Value = 0.385
$ws.Range("P20").Value = 1.7
$ws.Range("Q20").Value = 3.9
$ws.Range("R20").Value = 0.452
$ws.Range("S20").Value = 4.1
$ws.Range("T20").Value = 5.2
$ws.Range("U20").Value = 0.772
$ws.Range("Y20").Value = 5.6
$ws.Range("Z20").Value = 0.8
$ws.Range("AD20").Value = 20
$ws.Range("G22").Value = 307
$ws.Range("I22").Value = 7910
$ws.Range("T22").Value = 3.2
$ws.Range("U22").Value = 0.8149999999999999
$ws.Range("T24").Value = 1.6
$ws.Range("U24").Value = 0.06000000000000005

$ws = $wb.Worksheets.Item("per_poss")
$ws.Range("G11").Value = 67
$ws.Range("I11").Value = 1202
$ws.Range("J11").Value = 6
$ws.Range("K11").Value = 15.6
$ws.Range("L11").Value = 0.382
$ws.Range("N11").Value = 12.3
$ws.Range("R11").Value = 0.458
$ws.Range("S11").Value = 3.5
$ws.Range("T11").Value = 4.6
$ws.Range("U11").Value = 0.759
$ws.Range("V11").Value = 0.8
$ws.Range("W11").Value = 3.6
$ws.Range("Y11").Value = 7.1
$ws.Range("AD11").Value = 19.9
$ws.Range("AF11").Value = 120
$ws.Range("G13").Value = 14
$ws.Range("I13").Value = 391
$ws.Range("K13").Value = 19.3
$ws.Range("L13").Value = 0.403
$ws.Range("M13").Value = 5.5
$ws.Range("N13").Value = 14.2
$ws.Range("O13").Value = 0.385
$ws.Range("Q13").Value = 5.1
$ws.Range("R13").Value = 0.452
$ws.Range("S13").Value = 5.3
$ws.Range("T13").Value = 6.9
$ws.Range("U13").Value = 0.772
$ws.Range("W13").Value = 3.8
$ws.Range("X13").Value = 4.6
$ws.Range("Y13").Value = 7.4
$ws.Range("Z13").Value = 1.1
$ws.Range("AB13").Value = 2.3
$ws.Range("AC13").Value = 2.8
$ws.Range("AD13").Value = 26.3
$ws.Range("AF13").Value = 122
$ws.Range("AF17").Value = 111
$ws.Range("AF19").Value = 1
$ws.Range("G20").Value = 14
$ws.Range("I20").Value = 391
$ws.Range("K20").Value = 19.3
$ws.Range("L20").Value = 0.403
$ws.Range("M20").Value = 5.5
$ws.Range("N20").Value = 14.2
$ws.Range("O20").Value = 0.385
$ws.Range("Q20").Value = 5.1
$ws.Range("R20").Value = 0.452
$ws.Range("S20").Value = 5.3
$ws.Range("T20").Value = 6.9
$ws.Range("U20").Value = 0.772
$ws.Range("W20").Value = 3.8
$ws.Range("X20").Value = 4.6
$ws.Range("Y20").Value = 7.4
$ws.Range("Z20").Value = 1.1
$ws.Range("AB20").Value = 2.3
$ws.Range("AC20").Value = 2.8
$ws.Range("AD20").Value = 26.3
$ws.Range("AF20").Value = 122
$ws.Range("G22").Value = 307
$ws.Range("I22").Value = 7910
$ws.Range("N22").Value = 12.8
$ws.Range("U22").Value = 0.8149999999999999
$ws.Range("AF22").Value = 111
$ws.Range("N24").Value = -0.8000000000000007
$ws.Range("U24").Value = 0.06000000000000005
$ws.Range("AF24").Value = 1

$ws = $wb.Worksheets.Item("advanced")
$ws.Range("G11").Value = 67
$ws.Range("H11").Value = 1202
$ws.Range("I11").Value = 13.2
$ws.Range("J11").Value = 0.5629999999999999
$ws.Range("K11").Value = 0.787
$ws.Range("L11").Value = 0.297
$ws.Range("N11").Value = 8.6
$ws.Range("P11").Value = 19.7
$ws.Range("R11").Value = 1.1
$ws.Range("T11").Value = 17.1
$ws.Range("Y11").Value = 0.102
$ws.Range("AA11").Value = 0.7
$ws.Range("AC11").Value = 0.6
$ws.Range("Y12").Value = 0.101
$ws.Range("AB12").Value = 0.5
$ws.Range("G13").Value = 14
$ws.Range("H13").Value = 391
$ws.Range("I13").Value = 16
$ws.Range("J13").Value = 0.589
$ws.Range("K13").Value = 0.736
$ws.Range("L13").Value = 0.358
$ws.Range("M13").Value = 1.9
$ws.Range("N13").Value = 9.1
$ws.Range("O13").Value = 5.3
$ws.Range("P13").Value = 21.7
$ws.Range("Q13").Value = 1.1
$ws.Range("R13").Value = 1
$ws.Range("S13").Value = 9.4
$ws.Range("T13").Value = 21.5
$ws.Range("X13").Value = 0.8
$ws.Range("Y13").Value = 0.104
$ws.Range("AA13").Value = 3.1
$ws.Range("AB13").Value = -1.3
$ws.Range("AC13").Value = 1.8
$ws.Range("Y17").Value = 0.073
$ws.Range("AB17").Value = -0.7
$ws.Range("AC17").Value = -0.7
$ws.Range("Y19").Value = -0.02399999999999999
$ws.Range("AB19").Value = -0.3
$ws.Range("AC19").Value = -0.4000000000000001
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 391
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 0.589
$ws.Range("K20").Value = 0.736
$ws.Range("L20").Value = 0.358
$ws.Range("M20").Value = 1.9
$ws.Range("N20").Value = 9.1
$ws.Range("O20").Value = 5.3
$ws.Range("P20").Value = 21.7
$ws.Range("Q20").Value = 1.1
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 9.4
$ws.Range("T20").Value = 21.5
$ws.Range("X20").Value = 0.8
$ws.Range("Y20").Value = 0.104
$ws.Range("AA20").Value = 3.1
$ws.Range("AB20").Value = -1.3
$ws.Range("AC20").Value = 1.8
$ws.Range("G22").Value = 307
$ws.Range("H22").Value = 7910
$ws.Range("L22").Value = 0.227
$ws.Range("N22").Value = 8.199999999999999
$ws.Range("P22").Value = 26.2
$ws.Range("V22").Value = 8.300000000000001
$ws.Range("AC22").Value = -0.3
$ws.Range("AD22").Value = 3.4
$ws.Range("L24").Value = 0.217
$ws.Range("N24").Value = 6.5
$ws.Range("P24").Value = -16.8
$ws.Range("V24").Value = -8.200000000000001
$ws.Range("AC24").Value = -0.8
$ws.Range("AD24").Value = -3.4
